$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2862.85
$ws.Range("J17").Value = 2960.8948
$ws.Range("L17").Value = 8882.6844
$ws.Range("N17").Value = -9218.6844

$ws.Range("H132").Value = 30984.705
$ws.Range("I132").Value = 36992.18
$ws.Range("K132").Value = 110976.54
$ws.Range("M132").Value = -108446.54

$ws.Range("H137").Value = 7694172.5
$ws.Range("I137").Value = 1602
$ws.Range("K137").Value = 4806
$ws.Range("M137").Value = -2256

$ws.Range("H138").Value = 14155.125
$ws.Range("J138").Value = 8153.778
$ws.Range("L138").Value = 24461.334
$ws.Range("N138").Value = -34741.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1076.0952
$ws.Range("I2").Value = 611.1539
$ws.Range("K2").Value = 611.1539
$ws.Range("M2").Value = -498.1539

$ws.Range("H32").Value = 1447.8125
$ws.Range("I32").Value = 943.0339
$ws.Range("K32").Value = 943.0339
$ws.Range("M32").Value = -656.0339

$ws.Range("H45").Value = 47212.547
$ws.Range("I45").Value = 60463.47
$ws.Range("J45").Value = 2159.4
$ws.Range("K45").Value = 60463.47
$ws.Range("L45").Value = 2159.4
$ws.Range("M45").Value = -60086.47
$ws.Range("N45").Value = -2913.4

$ws.Range("H61").Value = 1438743.6
$ws.Range("I61").Value = 3721.6333
$ws.Range("J61").Value = 5352440
$ws.Range("K61").Value = 3721.6333
$ws.Range("L61").Value = 5352440
$ws.Range("M61").Value = -3509.6333
$ws.Range("N61").Value = -5352864

$ws.Range("H74").Value = 414343.25
$ws.Range("I74").Value = 1041.079
$ws.Range("J74").Value = 1395935.9
$ws.Range("K74").Value = 1041.079
$ws.Range("L74").Value = 1395935.9
$ws.Range("M74").Value = -167.079
$ws.Range("N74").Value = -1397683.9

$ws.Range("H77").Value = 414343.25
$ws.Range("I77").Value = 1041.079
$ws.Range("J77").Value = 1395935.9
$ws.Range("K77").Value = 5205.395
$ws.Range("L77").Value = 6979679.5
$ws.Range("M77").Value = -837.3949999999995
$ws.Range("N77").Value = -6988415.5

$ws.Range("H97").Value = 3463.611
$ws.Range("I97").Value = 4673.32
$ws.Range("K97").Value = 4673.32
$ws.Range("M97").Value = -4177.32

$ws.Range("H102").Value = 1083.4849
$ws.Range("I102").Value = 890.931
$ws.Range("K102").Value = 890.931
$ws.Range("M102").Value = 731.069

$ws.Range("H110").Value = 2741.8572
$ws.Range("I110").Value = 2698.8333
$ws.Range("K110").Value = 2698.8333
$ws.Range("M110").Value = -653.8332999999998

$ws.Range("H116").Value = 1076.0952
$ws.Range("I116").Value = 611.1539
$ws.Range("K116").Value = 611.1539
$ws.Range("M116").Value = 1682.8461

$ws.Range("H122").Value = 1820.5918
$ws.Range("I122").Value = 1623.6666
$ws.Range("K122").Value = 4870.9998
$ws.Range("M122").Value = -2420.9998

$ws.Range("H136").Value = 1438743.6
$ws.Range("I136").Value = 3721.6333
$ws.Range("J136").Value = 5352440
$ws.Range("K136").Value = 11164.8999
$ws.Range("L136").Value = 16057320
$ws.Range("M136").Value = -8614.8999
$ws.Range("N136").Value = -16062420

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1076.0952
$ws.Range("I3").Value = 611.1539
$ws.Range("K3").Value = 611.1539
$ws.Range("M3").Value = -497.1539

$ws.Range("H35").Value = 69995
$ws.Range("J35").Value = 69995
$ws.Range("L35").Value = 69995
$ws.Range("N35").Value = -70615

$ws.Range("H86").Value = 2003.4286
$ws.Range("I86").Value = 746.86664
$ws.Range("K86").Value = 746.86664
$ws.Range("M86").Value = 376.13336

$ws.Range("H89").Value = 2003.4286
$ws.Range("I89").Value = 746.86664
$ws.Range("K89").Value = 3734.3332
$ws.Range("M89").Value = 1881.6668

$ws.Range("H94").Value = 1014.9655
$ws.Range("I94").Value = 1137.44
$ws.Range("J94").Value = 249.5
$ws.Range("K94").Value = 1137.44
$ws.Range("L94").Value = 249.5
$ws.Range("M94").Value = -686.4400000000001
$ws.Range("N94").Value = -1151.5

$ws.Range("H105").Value = 6069
$ws.Range("I105").Value = 7326.4375
$ws.Range("J105").Value = 3833.5557
$ws.Range("K105").Value = 7326.4375
$ws.Range("L105").Value = 3833.5557
$ws.Range("M105").Value = -5579.4375
$ws.Range("N105").Value = -7327.5557

$ws.Range("H107").Value = 7453.1914
$ws.Range("I107").Value = 8222.462
$ws.Range("J107").Value = 3703
$ws.Range("K107").Value = 8222.462
$ws.Range("L107").Value = 3703
$ws.Range("M107").Value = -6302.462
$ws.Range("N107").Value = -7543

$ws.Range("H134").Value = 12501547
$ws.Range("I134").Value = 1402.8524
$ws.Range("K134").Value = 4208.5572
$ws.Range("M134").Value = -1673.5572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10002500
$ws.Range("I4").Value = 5000
$ws.Range("K4").Value = 5000
$ws.Range("M4").Value = -4888

$ws.Range("H5").Value = 1830.3334
$ws.Range("I5").Value = 3404.3333
$ws.Range("J5").Value = 256.33334
$ws.Range("K5").Value = 3404.3333
$ws.Range("L5").Value = 256.33334
$ws.Range("M5").Value = -3292.3333
$ws.Range("N5").Value = -480.33334

$ws.Range("H31").Value = 3350.3845
$ws.Range("I31").Value = 2832.4333
$ws.Range("K31").Value = 2832.4333
$ws.Range("M31").Value = -2537.4333

$ws.Range("H34").Value = 3350.3845
$ws.Range("I34").Value = 2832.4333
$ws.Range("K34").Value = 2832.4333
$ws.Range("M34").Value = -2630.4333

$ws.Range("H58").Value = 2082.7437
$ws.Range("I58").Value = 1848.5862
$ws.Range("K58").Value = 1848.5862
$ws.Range("M58").Value = -1645.5862

$ws.Range("H105").Value = 2154.5557
$ws.Range("I105").Value = 1713.7142
$ws.Range("J105").Value = 3697.5
$ws.Range("K105").Value = 1713.7142
$ws.Range("L105").Value = 3697.5
$ws.Range("M105").Value = 33.28580000000011
$ws.Range("N105").Value = -7191.5

$ws.Range("H132").Value = 62066.35
$ws.Range("I132").Value = 69675.2
$ws.Range("K132").Value = 209025.6
$ws.Range("M132").Value = -206495.6

$ws.Range("H134").Value = 1804.8823
$ws.Range("I134").Value = 1311
$ws.Range("K134").Value = 3933
$ws.Range("M134").Value = -1398

$ws.Range("H136").Value = 2082.7437
$ws.Range("I136").Value = 1848.5862
$ws.Range("K136").Value = 5545.7586
$ws.Range("M136").Value = -2995.7586

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1157.25
$ws.Range("J2").Value = 1629.9546
$ws.Range("L2").Value = 9779.7276
$ws.Range("N2").Value = -10005.7276

$ws.Range("H8").Value = 2537.25
$ws.Range("I8").Value = 2537.25
$ws.Range("K8").Value = 7611.75
$ws.Range("M8").Value = -7472.75

$ws.Range("H23").Value = 143892.72
$ws.Range("I23").Value = 3521.5
$ws.Range("J23").Value = 200041.2
$ws.Range("K23").Value = 10564.5
$ws.Range("L23").Value = 600123.6000000001
$ws.Range("M23").Value = -10329.5
$ws.Range("N23").Value = -600593.6000000001

$ws.Range("H38").Value = 1123
$ws.Range("I38").Value = 565.75
$ws.Range("K38").Value = 1697.25
$ws.Range("M38").Value = -1350.25

$ws.Range("H55").Value = 4699.2383
$ws.Range("I55").Value = 1449.8334
$ws.Range("K55").Value = 4349.5002
$ws.Range("M55").Value = -4172.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 25497.25
$ws.Range("J98").Value = 25497.25
$ws.Range("L98").Value = 25497.25
$ws.Range("N98").Value = -31487.25

$ws.Range("H113").Value = 2148
$ws.Range("I113").Value = 1904.7142
$ws.Range("K113").Value = 1904.7142
$ws.Range("M113").Value = 265.2858000000001

$ws.Range("H122").Value = 2177.7273
$ws.Range("I122").Value = 993
$ws.Range("K122").Value = 2979
$ws.Range("M122").Value = -529

$ws.Range("H132").Value = 957188.25
$ws.Range("I132").Value = 10327.857
$ws.Range("J132").Value = 1976884.1
$ws.Range("K132").Value = 30983.571
$ws.Range("L132").Value = 5930652.300000001
$ws.Range("M132").Value = -28453.571
$ws.Range("N132").Value = -5935712.300000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10510.368
$ws.Range("I7").Value = 6840.5454
$ws.Range("K7").Value = 6840.5454
$ws.Range("M7").Value = -6728.5454

$ws.Range("H46").Value = 14888.777
$ws.Range("J46").Value = 7000
$ws.Range("L46").Value = 7000
$ws.Range("N46").Value = -7376

$ws.Range("H55").Value = 840.5
$ws.Range("I55").Value = 385.25
$ws.Range("J55").Value = 1068.125
$ws.Range("K55").Value = 385.25
$ws.Range("L55").Value = 1068.125
$ws.Range("M55").Value = -212.25
$ws.Range("N55").Value = -1414.125

$ws.Range("H126").Value = 10510.368
$ws.Range("I126").Value = 6840.5454
$ws.Range("K126").Value = 20521.6362
$ws.Range("M126").Value = -18051.6362

$ws.Range("H132").Value = 3828.3235
$ws.Range("I132").Value = 3773.3
$ws.Range("K132").Value = 11319.9
$ws.Range("M132").Value = -8789.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 22674672
$ws.Range("I81").Value = 22674672
$ws.Range("K81").Value = 45349344
$ws.Range("M81").Value = -45348283

$ws.Range("H84").Value = 22674672
$ws.Range("I84").Value = 22674672
$ws.Range("K84").Value = 226746720
$ws.Range("M84").Value = -226741416
